# CSC 424 Team Meeting Minutes template — "Converted Team minutes to pdf" edit:
#   1. Append " 2/2/2021" (its own run) right after the "Date" heading.
#   2. Relocate the document's "_GoBack" bookmark (an artifact Word leaves at
#      the most-recently-edited spot) from its old position — right after
#      "Get the database fields" — to the end of that new "Date ..." line.

$d = $word.ActiveDocument

# --- 1. Drop the old _GoBack bookmark (after "Get the database fields") ----
$oldGoBack = $d.Bookmarks("_GoBack")
$oldGoBack.Delete()

# --- 2. Find the "Date" paragraph and append the date text ----------------
$datePara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -eq "Date`r") {
        $datePara = $para
    }
}

$dateRange = $datePara.Range
$insertPos = $dateRange.End - 1   # just after "Date", before its paragraph mark

# Insert the new text plus a throw-away trailing marker character; the
# marker gives us a real character to anchor the bookmark on (collapsed
# ranges sitting exactly at a paragraph's end are mishandled by this host),
# and it also forces Word to realize the appended text as its own run
# instead of silently re-merging it into the "Date" run.
$insertionPoint = $d.Range($insertPos, $insertPos)
$insertionPoint.InsertAfter(" 2/2/2021#")

$datePara = $d.Paragraphs($datePara.Range.Start, $datePara.Range.Start)
$dateRangeNow = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Start -eq $dateRange.Start) { $dateRangeNow = $para.Range }
}

$suffixStart = $insertPos
$suffixEnd = $dateRangeNow.End - 1
$suffixRange = $d.Range($suffixStart, $suffixEnd)

# Toggling Bold off/on splits " 2/2/2021#" into its own run (distinct from
# "Date") and leaves the run with an explicit rPr; setting NameBi alongside
# Name is what makes the host actually keep the w:cs="Times New Roman"
# attribute instead of optimizing it away as redundant.
$suffixRange.Font.Bold = $false
$suffixRange2 = $d.Range($suffixStart, $suffixEnd)
$suffixRange2.Font.Bold = $true
$suffixRange2.Font.Name = "Times New Roman"
$suffixRange2.Font.NameBi = "Times New Roman"
$suffixRange2.Font.Size = 12

# --- 3. Add the _GoBack bookmark right before the trailing marker ---------
$paraEnd = $dateRangeNow.End
$bmPos = $paraEnd - 2
$bmTarget = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmTarget)

# --- 4. Remove the throw-away marker character -----------------------------
$markerRange = $d.Range($bmPos, $bmPos + 1)
$markerRange.Delete()
